# Weekly update: insert a new "Poroto granado" price record as the new
# most-recent row (row 91), pushing the existing historical rows down by
# one (old row 91 -> 92, ... old row 119 -> 120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91:119 down to 92:120, leaving a blank row 91 for the new record.
$ws.Rows("91:91").Insert()

# Populate the new record in row 91.
$ws.Range("A91").Value = 7
$ws.Range("B91").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C91").Value = "Ñuble"
$ws.Range("D91").Value = 44985
$ws.Range("E91").Value = 16
$ws.Range("F91").Value = 100112030
$ws.Range("G91").Value = "Poroto granado"
$ws.Range("H91").Value = "Sin especificar"
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 60
$ws.Range("K91").Value = 26000
$ws.Range("L91").Value = 27000
$ws.Range("M91").Value = 26500
$ws.Range("N91").Value = "`$/saco 25 kilos"
$ws.Range("O91").Value = "Provincia de Diguillín"
$ws.Range("P91").Value = 1060
$ws.Range("Q91").Value = 25
$ws.Range("R91").Value = "Hortaliza"
